$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.047.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "'2.419.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'552.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "

# Row 6
$ws.Range("D6").Value = "'137.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.57%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +3.66%  "

# Row 9
$ws.Range("E9").Value = "  -1.58%  "

# Row 10
$ws.Range("E10").Value = "  -2.48%  "

# Row 11
$ws.Range("E11").Value = "  -0.93%  "

# Row 12
$ws.Range("E12").Value = "  -2.33%  "

# Row 13
$ws.Range("E13").Value = "  +0.44%  "

# Row 14
$ws.Range("D14").Value = "'2.850.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "

# Row 15
$ws.Range("D15").Value = "'59.989.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "

# Row 16
$ws.Range("E16").Value = "  -1.90%  "

# Row 17
$ws.Range("D17").Value = "'2.417.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18
$ws.Range("D18").Value = "'11.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

# Row 19
$ws.Range("E19").Value = "  -0.43%  "

# Row 20
$ws.Range("D20").Value = "'327.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "

# Row 21
$ws.Range("D21").Value = "'6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.71%  "

# Row 22
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").Value = "'65.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "

# Row 24
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
$ws.Range("E25").Value = "  +0.65%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("D27").Value = "'1.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "

# Row 28
$ws.Range("E28").Value = "  -2.85%  "

# Row 29
$ws.Range("E29").Value = "  -2.44%  "

# Row 30
$ws.Range("D30").Value = "'168.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "

# Row 31
$ws.Range("E31").Value = "  -4.39%  "

# Row 32
$ws.Range("E32").Value = "  +1.27%  "

# Row 33
$ws.Range("E33").Value = "  -1.47%  "

# Row 34
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("E35").Value = "  -0.78%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  -2.32%  "

# Row 38
$ws.Range("E38").Value = "  -2.27%  "

# Row 39
$ws.Range("D39").Value = "'324.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40
$ws.Range("D40").Value = "'0.404"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.48%  "

# Row 41
$ws.Range("D41").Value = "'3.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.15%  "

# Row 42
$ws.Range("D42").Value = "'140.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.83%  "

# Row 44
$ws.Range("E44").Value = "  -1.98%  "

# Row 45
$ws.Range("D45").Value = "'0.0517"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.77%  "

# Row 46
$ws.Range("E46").Value = "  +0.33%  "

# Row 47
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").Value = "'0.398"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.82%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0223"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49
$ws.Range("E49").Value = "  -0.09%  "

# Row 50
$ws.Range("E50").Value = "  -4.76%  "

# Row 51
$ws.Range("E51").Value = "  -1.05%  "
